$wb = $excel.ActiveWorkbook

# Duplicate the "Turkey" sheet (closest existing template) to the end of the
# workbook, producing the new "Croatia" market sheet.
$turkey = $wb.Worksheets.Item("Turkey")
$turkey.Activate() | Out-Null
$turkey.Cells.Select() | Out-Null
$turkey.Copy($null, $turkey) | Out-Null
$croatia = $wb.Worksheets.Item($wb.Worksheets.Count)
$croatia.Name = "Croatia"

# Fill in the market-specific values for the new sheet.
$croatia.Range("B2").Value = "Croatia Market"
$croatia.Range("B4").Value = "NGC-3139/T2424/T2415"

# Resize column B to fit the new (longer) content, matching the width Excel
# computed for the new market code.
$croatia.Columns("B").ColumnWidth = 24

# Leave the cursor on B4, matching the authored selection state.
$croatia.Range("B4").Select() | Out-Null
